$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value could be misread as a number by Excel-style
# auto-detection; force Text format so they stay plain strings, then
# clear the format override so no stray style survives on the cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "274.94"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5046"
$ws.Range("D7").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.45"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06641"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.96"
$ws.Range("D11").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07806"
$ws.Range("D13").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.052"
$ws.Range("D15").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007988"
$ws.Range("D20").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.721"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.959"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.050"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.64"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.132"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.94"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "108.44"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.327"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.191"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08761"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04798"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7236"
$ws.Range("D34").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.878"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9990"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.030"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01855"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5168"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.272"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9456"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.23"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.165"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.973"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9992"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1377"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4548"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.299"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.96"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.494"
$ws.Range("D51").ClearFormats()

# Remaining cells are unambiguous text already (contain letters, %,
# extra separators, or URLs) so a plain Value assignment is safe.
$ws.Range("D2").Value = "25.738.72"
$ws.Range("E2").Value = "  -5.74%  "
$ws.Range("D3").Value = "1.808.86"
$ws.Range("E3").Value = "  -5.04%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -10.28%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -7.06%  "
$ws.Range("E8").Value = "  -7.92%  "
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("E10").Value = "  -9.12%  "
$ws.Range("E11").Value = "  -9.87%  "
$ws.Range("E12").Value = "  -7.78%  "
$ws.Range("E13").Value = "  -4.60%  "
$ws.Range("D14").Value = "1.801.13"
$ws.Range("E14").Value = "  +34.24%  "
$ws.Range("E15").Value = "  -5.64%  "
$ws.Range("E16").Value = "  -8.73%  "
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  -6.78%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -7.80%  "
$ws.Range("D21").Value = "25.777.57"
$ws.Range("E21").Value = "  -5.63%  "
$ws.Range("E22").Value = "  -6.61%  "
$ws.Range("E23").Value = "  -8.07%  "
$ws.Range("E24").Value = "  -7.16%  "
$ws.Range("E25").Value = "  -4.71%  "
$ws.Range("E26").Value = "  -7.86%  "
$ws.Range("E27").Value = "  -5.69%  "
$ws.Range("E28").Value = "  -7.72%  "
$ws.Range("E29").Value = "  -6.98%  "
$ws.Range("E30").Value = "  -10.85%  "
$ws.Range("E31").Value = "  -10.42%  "
$ws.Range("E32").Value = "  -4.70%  "
$ws.Range("E34").Value = "  -12.54%  "
$ws.Range("E35").Value = "  -7.64%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("E38").Value = "  -8.70%  "
$ws.Range("E39").Value = "  -7.15%  "
$ws.Range("E40").Value = "  -13.98%  "
$ws.Range("E41").Value = "  -15.78%  "
$ws.Range("E42").Value = "  -12.15%  "
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("E45").Value = "  -14.15%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("E47").Value = "  -10.04%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("E48").Value = "  -11.73%  "
$ws.Range("E49").Value = "  -8.49%  "
$ws.Range("E50").Value = "  -5.68%  "
$ws.Range("E51").Value = "  -8.81%  "
